# Weekly refresh of the "Rabanito" (radish) wholesale price series.
#
# A new week's observation (2022-10-26, Volumen 7900) is inserted as the
# most recent record right after the header row, pushing every existing
# record down by one row (dimension grows from A1:R32 to A1:R33). All of
# the "moving" columns for a given record (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) travel together with
# it; the columns that are constant for this whole sheet (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Variedad, Unidad de
# comercializacion, Kg o Unidades, Clasificacion) are simply re-stamped on
# every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot row 2 (the record that is about to be pushed down to row 3) so
# its values can be re-written into the freshly inserted row 2.
$A2 = $ws.Range("A2").Value2
$B2 = $ws.Range("B2").Value2
$C2 = $ws.Range("C2").Value2
$E2 = $ws.Range("E2").Value2
$F2 = $ws.Range("F2").Value2
$G2 = $ws.Range("G2").Value2
$H2 = $ws.Range("H2").Value2
$I2 = $ws.Range("I2").Value2
$K2 = $ws.Range("K2").Value2
$L2 = $ws.Range("L2").Value2
$M2 = $ws.Range("M2").Value2
$N2 = $ws.Range("N2").Value2
$O2 = $ws.Range("O2").Value2
$P2 = $ws.Range("P2").Value2
$Q2 = $ws.Range("Q2").Value2
$R2 = $ws.Range("R2").Value2

# Insert a brand-new blank row above the current row 2; this shifts rows
# 2-32 down to 3-33 and grows the sheet dimension to A1:R33.
$ws.Rows(2).Insert()

# The inserted row inherits stray formatting from the insert op; start it
# clean so the rebuilt cells pick up the sheet's normal (unstyled) look.
$ws.Rows(2).ClearFormats()

# Populate the new row 2 with the latest observation: same dimensions as
# the old row 2 had, but a new date and a new Volumen.
$ws.Range("A2").Value = $A2
$ws.Range("B2").Value = $B2
$ws.Range("C2").Value = $C2
$ws.Range("D2").Value = 44860
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = $E2
$ws.Range("F2").Value = $F2
$ws.Range("G2").Value = $G2
$ws.Range("H2").Value = $H2
$ws.Range("I2").Value = $I2
$ws.Range("J2").Value = 7900
$ws.Range("K2").Value = $K2
$ws.Range("L2").Value = $L2
$ws.Range("M2").Value = $M2
$ws.Range("N2").Value = $N2
$ws.Range("O2").Value = $O2
$ws.Range("P2").Value = $P2
$ws.Range("Q2").Value = $Q2
$ws.Range("R2").Value = $R2
